# Applies the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force the cell to stay a text cell even when the new value looks like
    # a number (e.g. "213.10"), matching the inlineStr-typed source data,
    # then drop back to the default "Normal" style so no stray number format
    # is left attached to the cell.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "27.682.42"
$ws.Range("E2").Value = "  +0.98%  "

# Row 3
$ws.Range("D3").Value = "1.644.11"
$ws.Range("E3").Value = "  +0.22%  "

# Row 5
Set-TextValue "D5" "213.10"
$ws.Range("E5").Value = "  +0.75%  "

# Row 6
$ws.Range("E6").Value = "  -0.59%  "

# Row 7
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
Set-TextValue "D8" "23.30"
$ws.Range("E8").Value = "  +0.94%  "

# Row 9
$ws.Range("E9").Value = "  +1.20%  "

# Row 10
$ws.Range("E10").Value = "  +0.72%  "

# Row 11
Set-TextValue "D11" "0.0895"
$ws.Range("E11").Value = "  +0.51%  "

# Row 12
$ws.Range("D12").Value = "1.877.05"

# Row 13
$ws.Range("D13").Value = "1.644.46"
$ws.Range("E13").Value = "  +0.03%  "

# Row 14
$ws.Range("E14").Value = "  +0.37%  "

# Row 15
$ws.Range("E15").Value = "  +1.19%  "

# Row 16
Set-TextValue "D16" "64.67"
$ws.Range("E16").Value = "  +0.79%  "

# Row 17
$ws.Range("D17").Value = "27.662.73"
$ws.Range("E17").Value = "  +1.02%  "

# Row 18
Set-TextValue "D18" "230.92"
$ws.Range("E18").Value = "  +0.98%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0723"
$ws.Range("E19").Value = "  +0.73%  "

# Row 20
Set-TextValue "D20" "7.64"

# Row 21
$ws.Range("E21").Value = "  +0.10%  "

# Row 22
$ws.Range("E22").Value = "  -0.57%  "

# Row 23
Set-TextValue "D23" "10.03"
$ws.Range("E23").Value = "  +8.11%  "

# Row 24
$ws.Range("E24").Value = "  -2.78%  "

# Row 25
Set-TextValue "D25" "149.91"
$ws.Range("E25").Value = "  +1.45%  "

# Row 26
Set-TextValue "D26" "6.93"
$ws.Range("E26").Value = "  +0.11%  "

# Row 27
$ws.Range("E27").Value = "  -1.12%  "

# Row 28
$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D28" "1.00"
$ws.Range("E28").Value = "  -0.04%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D29" "15.63"
$ws.Range("E29").Value = "  +0.83%  "

# Row 30
$ws.Range("E30").Value = "  +0.37%  "

# Row 31
$ws.Range("E31").Value = "  +0.85%  "

# Row 32
$ws.Range("E32").Value = "  +0.66%  "

# Row 33
$ws.Range("D33").Value = "1.444.23"
$ws.Range("E33").Value = "  +2.54%  "

# Row 34
$ws.Range("E34").Value = "  +1.38%  "

# Row 35
$ws.Range("E35").Value = "  +1.65%  "

# Row 36
$ws.Range("E36").Value = "  -1.31%  "

# Row 37
Set-TextValue "D37" "0.569"
$ws.Range("E37").Value = "  +1.29%  "

# Row 38
Set-TextValue "D38" "0.881"
$ws.Range("E38").Value = "  +0.41%  "

# Row 39
Set-TextValue "D39" "0.0167"
$ws.Range("E39").Value = "  +0.32%  "

# Row 40
Set-TextValue "D40" "0.888"
$ws.Range("E40").Value = "  +12.38%  "

# Row 41
$ws.Range("E41").Value = "  -0.04%  "

# Row 42
$ws.Range("E42").Value = "  +0.08%  "

# Row 43
Set-TextValue "D43" "5.62"
$ws.Range("E43").Value = "  +3.16%  "

# Row 44
Set-TextValue "D44" "67.11"
$ws.Range("E44").Value = "  +4.21%  "

# Row 45
Set-TextValue "D45" "2.25"
$ws.Range("E45").Value = "  +1.17%  "

# Row 46
$ws.Range("D46").Value = "1.786.60"
$ws.Range("E46").Value = "  +0.17%  "

# Row 47
Set-TextValue "D47" "1.72"
$ws.Range("E47").Value = "  +5.35%  "

# Row 48
Set-TextValue "D48" "85.65"
$ws.Range("E48").Value = "  -1.78%  "

# Row 49
$ws.Range("E49").Value = "  +0.36%  "

# Row 50
Set-TextValue "D50" "7.74"
$ws.Range("E50").Value = "  +1.15%  "

# Row 51
$ws.Range("E51").Value = "  +0.92%  "
